$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.853.05"
$ws.Range("E2").Value = "  +2.76%  "

$ws.Range("D3").Value = "2.955.15"
$ws.Range("E3").Value = "  +0.94%  "

$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.999"
$cell.ClearFormats()
$ws.Range("E4").Value = "  -0.08%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "595.17"
$cell.ClearFormats()
$ws.Range("E5").Value = "  +0.37%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "146.15"
$cell.ClearFormats()
$ws.Range("E6").Value = "  +1.21%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.ClearFormats()
$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("D8").Value = "2.957.75"
$ws.Range("E8").Value = "  +1.02%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.508"
$cell.ClearFormats()
$ws.Range("E9").Value = "  +0.88%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "7.26"
$cell.ClearFormats()
$ws.Range("E10").Value = "  +3.44%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.151"
$cell.ClearFormats()
$ws.Range("E11").Value = "  +5.92%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.445"
$cell.ClearFormats()
$ws.Range("E12").Value = "  +0.87%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000239"
$cell.ClearFormats()
$ws.Range("E13").Value = "  +6.10%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "33.10"
$cell.ClearFormats()
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("E15").Value = "  -0.47%  "

$ws.Range("D16").Value = "3.444.09"
$ws.Range("E16").Value = "  +0.85%  "

$ws.Range("D17").Value = "62.749.39"
$ws.Range("E17").Value = "  +2.69%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "6.74"
$cell.ClearFormats()
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "2.979.01"
$ws.Range("E19").Value = "  +1.63%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "442.80"
$cell.ClearFormats()
$ws.Range("E20").Value = "  +2.47%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "13.54"
$cell.ClearFormats()
$ws.Range("E21").Value = "  +0.46%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.669"
$cell.ClearFormats()
$ws.Range("E22").Value = "  -1.41%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.09"
$cell.ClearFormats()
$ws.Range("E23").Value = "  -0.11%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "81.43"
$cell.ClearFormats()
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "11.19"
$cell.ClearFormats()
$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "11.90"
$cell.ClearFormats()
$ws.Range("E26").Value = "  +0.95%  "

$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "2.13"
$cell.ClearFormats()
$ws.Range("E27").Value = "  -3.24%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.ClearFormats()

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "7.30"
$cell.ClearFormats()
$ws.Range("E29").Value = "  +4.79%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "2.62"
$cell.ClearFormats()
$ws.Range("E30").Value = "  +0.60%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "2.17"
$cell.ClearFormats()
$ws.Range("E31").Value = "  -2.95%  "

$ws.Range("D32").Value = "0.0₃0974"
$ws.Range("E32").Value = "  +10.05%  "

$ws.Range("E33").Value = "  -0.98%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "26.54"
$cell.ClearFormats()
$ws.Range("E34").Value = "  -0.80%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.994"
$cell.ClearFormats()
$ws.Range("E36").Value = "  -1.84%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.66"
$cell.ClearFormats()
$ws.Range("E37").Value = "  +0.38%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "3.10"
$cell.ClearFormats()
$ws.Range("E38").Value = "  +4.01%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.05"
$cell.ClearFormats()
$ws.Range("E39").Value = "  +2.52%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "49.59"
$cell.ClearFormats()
$ws.Range("E40").Value = "  -0.90%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "8.53"
$cell.ClearFormats()
$ws.Range("E41").Value = "  -1.09%  "

$ws.Range("E42").Value = "  -4.54%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.281"
$cell.ClearFormats()
$ws.Range("E43").Value = "  -0.52%  "

$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "40.68"
$cell.ClearFormats()
$ws.Range("E44").Value = "  -3.15%  "

$ws.Range("D45").Value = "2.726.14"
$ws.Range("E45").Value = "  +1.17%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "134.19"
$cell.ClearFormats()
$ws.Range("E46").Value = "  +0.60%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.0338"
$cell.ClearFormats()
$ws.Range("E47").Value = "  -2.80%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "362.84"
$cell.ClearFormats()
$ws.Range("E48").Value = "  -1.71%  "

$ws.Range("E49").Value = "  +0.06%  "

$ws.Range("E50").Value = "  -0.36%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "22.92"
$cell.ClearFormats()
$ws.Range("E51").Value = "  -3.70%  "
